$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "26.881.38"
$ws.Range("E2").Value = "  +0.87%  "
$ws.Range("D3").Value = "1.649.47"
$ws.Range("E3").Value = "  +0.40%  "
$ws.Range("E4").Value = "  +0.59%  "
$ws.Range("D5").Value = "'217.28"
$ws.Range("E5").Value = "  +0.92%  "
$ws.Range("E6").Value = "  -0.57%  "
$ws.Range("E7").Value = "  +0.61%  "
$ws.Range("D8").Value = "'0.253"
$ws.Range("E8").Value = "  +0.09%  "
$ws.Range("E9").Value = "  +0.22%  "
$ws.Range("E10").Value = "  -0.19%  "
$ws.Range("D11").Value = "'0.0845"
$ws.Range("E11").Value = "  +0.45%  "
$ws.Range("D12").Value = "1.664.37"
$ws.Range("E12").Value = "  +1.19%  "
$ws.Range("D13").Value = "'4.18"
$ws.Range("E13").Value = "  -0.40%  "
$ws.Range("E14").Value = "  -0.03%  "
$ws.Range("D15").Value = "'64.86"
$ws.Range("E15").Value = "  -0.93%  "
$ws.Range("D16").Value = "26.894.81"
$ws.Range("E16").Value = "  +0.76%  "
$ws.Range("E17").Value = "  -0.98%  "
$ws.Range("D18").Value = "'215.05"
$ws.Range("E18").Value = "  -1.11%  "
$ws.Range("E19").Value = "  +0.62%  "
$ws.Range("E20").Value = "  +1.08%  "
$ws.Range("E21").Value = "  +11.02%  "
$ws.Range("E22").Value = "  -0.61%  "
$ws.Range("E23").Value = "  -1.44%  "
$ws.Range("D24").Value = "'147.58"
$ws.Range("E24").Value = "  +1.46%  "
$ws.Range("E25").Value = "  +0.34%  "
$ws.Range("E26").Value = "  -0.95%  "
$ws.Range("D27").Value = "'7.19"
$ws.Range("E27").Value = "  +0.20%  "
$ws.Range("D28").Value = "'15.72"
$ws.Range("E28").Value = "  -0.36%  "
$ws.Range("D29").Value = "'0.0509"
$ws.Range("E29").Value = "  -1.71%  "
$ws.Range("D30").Value = "'1.19"
$ws.Range("E30").Value = "  +0.82%  "
$ws.Range("E31").Value = "  -0.50%  "
$ws.Range("E32").Value = "  -0.87%  "
$ws.Range("D33").Value = "1.297.33"
$ws.Range("E33").Value = "  +1.45%  "
$ws.Range("E34").Value = "  -0.20%  "
$ws.Range("E35").Value = "  +1.67%  "
$ws.Range("D36").Value = "'0.0178"
$ws.Range("E36").Value = "  -1.83%  "
$ws.Range("D37").Value = "'0.536"
$ws.Range("E37").Value = "  +0.89%  "
$ws.Range("E38").Value = "  -0.71%  "
$ws.Range("E39").Value = "  +0.65%  "
$ws.Range("E40").Value = "  -0.78%  "
$ws.Range("E41").Value = "  -0.52%  "
$ws.Range("D42").Value = "'5.35"
$ws.Range("E42").Value = "  -2.12%  "
$ws.Range("D43").Value = "1.787.86"
$ws.Range("E43").Value = "  +0.31%  "
$ws.Range("D44").Value = "'62.07"
$ws.Range("E44").Value = "  +3.75%  "
$ws.Range("D45").Value = "'92.11"
$ws.Range("E45").Value = "  +0.42%  "
$ws.Range("E46").Value = "  +1.25%  "
$ws.Range("D47").Value = "0.0₆0105"
$ws.Range("E47").Value = "  -0.55%  "
$ws.Range("E48").Value = "  +0.82%  "
$ws.Range("D49").Value = "'7.69"
$ws.Range("E49").Value = "  -1.68%  "
$ws.Range("D50").Value = "'0.0971"
$ws.Range("E50").Value = "  +0.04%  "
$ws.Range("E51").Value = "  +0.29%  "
